# Applies the commit "Add files via upload": clears several speaker-notes
# bodies down to an empty paragraph, rewrites the notes on slide 11 to a
# shorter remark, and merges a split "(10,1)" run back into one run on
# slide 2.

$p = $ppt.ActivePresentation
$rsquo = [char]0x2019

# ---------------------------------------------------------------------
# Slide notes that get fully cleared (every run removed, leaving just an
# empty paragraph in the "Notes Placeholder 2" shape).
# ---------------------------------------------------------------------
$notesToClear = @(1, 3, 4, 6, 7, 8, 9, 10)
foreach ($slideIdx in $notesToClear) {
    $slide = $p.Slides.Item($slideIdx)
    $notesShape = $slide.NotesPage.Shapes.Item(2)
    $notesShape.TextFrame.TextRange.Text = ""
}

# ---------------------------------------------------------------------
# Slide 11 notes: replace the long explanation with the shorter remark
# that starts mid-sentence ("these templates later, ...").
# ---------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$notes11 = $slide11.NotesPage.Shapes.Item(2)
$notes11.TextFrame.TextRange.Text = "these templates later, but for now let" + $rsquo + "s assume that the program has created templates, and let me explain how the program use these templates."

# ---------------------------------------------------------------------
# Slide 2: the "Biases" / "(10,1)" textbox had "(10,1" and ")" split
# across two runs; merge them back into a single run reading "(10,1)".
# The second paragraph already reads "(10,1)" so a direct re-assignment
# is a no-op for the engine's change-detection; force it through an
# intermediate value first.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shape = $slide2.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Biases(10,1)") {
            $biasesRange = $shape.TextFrame.TextRange
            $secondPara = $biasesRange.Paragraphs(2, 1)
            $secondPara.Text = "__tmp__"
            $freshRange = $shape.TextFrame.TextRange
            $secondParaFinal = $freshRange.Paragraphs(2, 1)
            $secondParaFinal.Text = "(10,1)"
        }
    }
}
